$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-09-21 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-22 Monday", 2) | Out-Null
$d.Content.Find.Execute("60-37=", $true, $false, $false, $false, $false, $true, 1, $false, "10+57=", 2) | Out-Null
$d.Content.Find.Execute("82-51=", $true, $false, $false, $false, $false, $true, 1, $false, "69-49=", 2) | Out-Null
$d.Content.Find.Execute("76+16=", $true, $false, $false, $false, $false, $true, 1, $false, "55+36=", 2) | Out-Null
$d.Content.Find.Execute("63-8=", $true, $false, $false, $false, $false, $true, 1, $false, "30+30=", 2) | Out-Null
$d.Content.Find.Execute("23-2=", $true, $false, $false, $false, $false, $true, 1, $false, "76+23=", 2) | Out-Null
$d.Content.Find.Execute("41-14=", $true, $false, $false, $false, $false, $true, 1, $false, "23+72=", 2) | Out-Null
$d.Content.Find.Execute("43+32=", $true, $false, $false, $false, $false, $true, 1, $false, "54-39=", 2) | Out-Null
$d.Content.Find.Execute("18+34=", $true, $false, $false, $false, $false, $true, 1, $false, "94-21=", 2) | Out-Null
$d.Content.Find.Execute("83-12=", $true, $false, $false, $false, $false, $true, 1, $false, "92-53=", 2) | Out-Null
$d.Content.Find.Execute("15+79=", $true, $false, $false, $false, $false, $true, 1, $false, "13+69=", 2) | Out-Null
$d.Content.Find.Execute("20-3=", $true, $false, $false, $false, $false, $true, 1, $false, "42+1=", 2) | Out-Null
$d.Content.Find.Execute("47+27=", $true, $false, $false, $false, $false, $true, 1, $false, "50+35=", 2) | Out-Null
$d.Content.Find.Execute("7+31=", $true, $false, $false, $false, $false, $true, 1, $false, "97-65=", 2) | Out-Null
$d.Content.Find.Execute("1+6=", $true, $false, $false, $false, $false, $true, 1, $false, "36-4=", 2) | Out-Null
$d.Content.Find.Execute("54+15=", $true, $false, $false, $false, $false, $true, 1, $false, "63-29=", 2) | Out-Null
$d.Content.Find.Execute("64+30=", $true, $false, $false, $false, $false, $true, 1, $false, "8+42=", 2) | Out-Null
$d.Content.Find.Execute("68-4=", $true, $false, $false, $false, $false, $true, 1, $false, "23-11=", 2) | Out-Null
$d.Content.Find.Execute("70-31=", $true, $false, $false, $false, $false, $true, 1, $false, "57-19=", 2) | Out-Null
$d.Content.Find.Execute("63+28=", $true, $false, $false, $false, $false, $true, 1, $false, "41+15=", 2) | Out-Null
$d.Content.Find.Execute("26+68=", $true, $false, $false, $false, $false, $true, 1, $false, "77-63=", 2) | Out-Null
$d.Content.Find.Execute("55-24=", $true, $false, $false, $false, $false, $true, 1, $false, "25+69=", 2) | Out-Null
$d.Content.Find.Execute("99-86=", $true, $false, $false, $false, $false, $true, 1, $false, "90-32=", 2) | Out-Null
$d.Content.Find.Execute("74+5=", $true, $false, $false, $false, $false, $true, 1, $false, "97-42=", 2) | Out-Null
$d.Content.Find.Execute("7+8=", $true, $false, $false, $false, $false, $true, 1, $false, "64+16=", 2) | Out-Null
$d.Content.Find.Execute("41-32=", $true, $false, $false, $false, $false, $true, 1, $false, "69-64=", 2) | Out-Null
$d.Content.Find.Execute("95-87=", $true, $false, $false, $false, $false, $true, 1, $false, "10+47=", 2) | Out-Null
$d.Content.Find.Execute("68-40=", $true, $false, $false, $false, $false, $true, 1, $false, "12+0=", 2) | Out-Null
$d.Content.Find.Execute("5+63=", $true, $false, $false, $false, $false, $true, 1, $false, "25+47=", 2) | Out-Null
$d.Content.Find.Execute("78-70=", $true, $false, $false, $false, $false, $true, 1, $false, "77-13=", 2) | Out-Null
$d.Content.Find.Execute("59-16=", $true, $false, $false, $false, $false, $true, 1, $false, "62-3=", 2) | Out-Null
$d.Content.Find.Execute("96-13=", $true, $false, $false, $false, $false, $true, 1, $false, "79-6=", 2) | Out-Null
$d.Content.Find.Execute("4+33=", $true, $false, $false, $false, $false, $true, 1, $false, "70-30=", 2) | Out-Null
$d.Content.Find.Execute("4+14=", $true, $false, $false, $false, $false, $true, 1, $false, "36-22=", 2) | Out-Null
$d.Content.Find.Execute("65-65=", $true, $false, $false, $false, $false, $true, 1, $false, "17+43=", 2) | Out-Null
$d.Content.Find.Execute("3+34=", $true, $false, $false, $false, $false, $true, 1, $false, "58-11=", 2) | Out-Null
$d.Content.Find.Execute("83-62=", $true, $false, $false, $false, $false, $true, 1, $false, "52-39=", 2) | Out-Null
$d.Content.Find.Execute("13+70=", $true, $false, $false, $false, $false, $true, 1, $false, "5+45=", 2) | Out-Null
$d.Content.Find.Execute("52-44=", $true, $false, $false, $false, $false, $true, 1, $false, "0+96=", 2) | Out-Null
$d.Content.Find.Execute("7+33=", $true, $false, $false, $false, $false, $true, 1, $false, "29-16=", 2) | Out-Null
$d.Content.Find.Execute("30-23=", $true, $false, $false, $false, $false, $true, 1, $false, "59-48=", 2) | Out-Null
$d.Content.Find.Execute("0+61=", $true, $false, $false, $false, $false, $true, 1, $false, "87-48=", 2) | Out-Null
$d.Content.Find.Execute("62-51=", $true, $false, $false, $false, $false, $true, 1, $false, "87-70=", 2) | Out-Null
$d.Content.Find.Execute("23+69=", $true, $false, $false, $false, $false, $true, 1, $false, "76-1=", 2) | Out-Null
$d.Content.Find.Execute("85-25=", $true, $false, $false, $false, $false, $true, 1, $false, "46+17=", 2) | Out-Null
$d.Content.Find.Execute("47+0=", $true, $false, $false, $false, $false, $true, 1, $false, "44+31=", 2) | Out-Null
$d.Content.Find.Execute("7+66=", $true, $false, $false, $false, $false, $true, 1, $false, "70+26=", 2) | Out-Null
$d.Content.Find.Execute("36+35=", $true, $false, $false, $false, $false, $true, 1, $false, "57+23=", 2) | Out-Null
$d.Content.Find.Execute("89-46=", $true, $false, $false, $false, $false, $true, 1, $false, "36+12=", 2) | Out-Null
$d.Content.Find.Execute("93-27=", $true, $false, $false, $false, $false, $true, 1, $false, "41+36=", 2) | Out-Null
$d.Content.Find.Execute("17+54=", $true, $false, $false, $false, $false, $true, 1, $false, "64-20=", 2) | Out-Null
$d.Content.Find.Execute("91-90=", $true, $false, $false, $false, $false, $true, 1, $false, "90-82=", 2) | Out-Null
$d.Content.Find.Execute("91+7=", $true, $false, $false, $false, $false, $true, 1, $false, "69-31=", 2) | Out-Null
$d.Content.Find.Execute("57-2=", $true, $false, $false, $false, $false, $true, 1, $false, "24+69=", 2) | Out-Null
$d.Content.Find.Execute("34+36=", $true, $false, $false, $false, $false, $true, 1, $false, "81-0=", 2) | Out-Null
$d.Content.Find.Execute("21+33=", $true, $false, $false, $false, $false, $true, 1, $false, "39+44=", 2) | Out-Null
$d.Content.Find.Execute("8+80=", $true, $false, $false, $false, $false, $true, 1, $false, "62-29=", 2) | Out-Null
$d.Content.Find.Execute("54+20=", $true, $false, $false, $false, $false, $true, 1, $false, "21+61=", 2) | Out-Null
$d.Content.Find.Execute("5+9=", $true, $false, $false, $false, $false, $true, 1, $false, "85-37=", 2) | Out-Null
$d.Content.Find.Execute("24+27=", $true, $false, $false, $false, $false, $true, 1, $false, "28-1=", 2) | Out-Null
$d.Content.Find.Execute("93-49=", $true, $false, $false, $false, $false, $true, 1, $false, "98-87=", 2) | Out-Null
$d.Content.Find.Execute("31-8=", $true, $false, $false, $false, $false, $true, 1, $false, "77-46=", 2) | Out-Null
$d.Content.Find.Execute("88-74=", $true, $false, $false, $false, $false, $true, 1, $false, "77+1=", 2) | Out-Null
$d.Content.Find.Execute("17-2=", $true, $false, $false, $false, $false, $true, 1, $false, "64-46=", 2) | Out-Null
$d.Content.Find.Execute("12+39=", $true, $false, $false, $false, $false, $true, 1, $false, "43+26=", 2) | Out-Null
$d.Content.Find.Execute("75-4=", $true, $false, $false, $false, $false, $true, 1, $false, "8+86=", 2) | Out-Null
$d.Content.Find.Execute("40-7=", $true, $false, $false, $false, $false, $true, 1, $false, "29+55=", 2) | Out-Null
$d.Content.Find.Execute("94-25=", $true, $false, $false, $false, $false, $true, 1, $false, "81+17=", 2) | Out-Null
$d.Content.Find.Execute("56-41=", $true, $false, $false, $false, $false, $true, 1, $false, "4+69=", 2) | Out-Null
$d.Content.Find.Execute("11+80=", $true, $false, $false, $false, $false, $true, 1, $false, "37+37=", 2) | Out-Null
$d.Content.Find.Execute("6+8=", $true, $false, $false, $false, $false, $true, 1, $false, "32+6=", 2) | Out-Null
$d.Content.Find.Execute("72-11=", $true, $false, $false, $false, $false, $true, 1, $false, "57+27=", 2) | Out-Null
$d.Content.Find.Execute("31+37=", $true, $false, $false, $false, $false, $true, 1, $false, "75-45=", 2) | Out-Null
$d.Content.Find.Execute("99-26=", $true, $false, $false, $false, $false, $true, 1, $false, "15+69=", 2) | Out-Null
$d.Content.Find.Execute("15-15=", $true, $false, $false, $false, $false, $true, 1, $false, "43-24=", 2) | Out-Null
$d.Content.Find.Execute("5+6=", $true, $false, $false, $false, $false, $true, 1, $false, "15+41=", 2) | Out-Null
$d.Content.Find.Execute("39+4=", $true, $false, $false, $false, $false, $true, 1, $false, "64-51=", 2) | Out-Null
$d.Content.Find.Execute("72-60=", $true, $false, $false, $false, $false, $true, 1, $false, "90+8=", 2) | Out-Null
$d.Content.Find.Execute("50+36=", $true, $false, $false, $false, $false, $true, 1, $false, "33+5=", 2) | Out-Null
$d.Content.Find.Execute("33+59=", $true, $false, $false, $false, $false, $true, 1, $false, "33+35=", 2) | Out-Null
$d.Content.Find.Execute("3+93=", $true, $false, $false, $false, $false, $true, 1, $false, "99-60=", 2) | Out-Null
$d.Content.Find.Execute("77-2=", $true, $false, $false, $false, $false, $true, 1, $false, "60-39=", 2) | Out-Null
$d.Content.Find.Execute("35-1=", $true, $false, $false, $false, $false, $true, 1, $false, "29+61=", 2) | Out-Null
$d.Content.Find.Execute("16+73=", $true, $false, $false, $false, $false, $true, 1, $false, "30+30=", 2) | Out-Null
$d.Content.Find.Execute("72+11=", $true, $false, $false, $false, $false, $true, 1, $false, "32+14=", 2) | Out-Null
$d.Content.Find.Execute("39+7=", $true, $false, $false, $false, $false, $true, 1, $false, "47-36=", 2) | Out-Null
$d.Content.Find.Execute("34-29=", $true, $false, $false, $false, $false, $true, 1, $false, "76+13=", 2) | Out-Null
$d.Content.Find.Execute("55+8=", $true, $false, $false, $false, $false, $true, 1, $false, "77-64=", 2) | Out-Null
$d.Content.Find.Execute("76-50=", $true, $false, $false, $false, $false, $true, 1, $false, "82-32=", 2) | Out-Null
$d.Content.Find.Execute("70-58=", $true, $false, $false, $false, $false, $true, 1, $false, "66+26=", 2) | Out-Null
$d.Content.Find.Execute("23-17=", $true, $false, $false, $false, $false, $true, 1, $false, "55+18=", 2) | Out-Null
$d.Content.Find.Execute("81-40=", $true, $false, $false, $false, $false, $true, 1, $false, "31+53=", 2) | Out-Null
$d.Content.Find.Execute("91-66=", $true, $false, $false, $false, $false, $true, 1, $false, "51+38=", 2) | Out-Null
$d.Content.Find.Execute("72+10=", $true, $false, $false, $false, $false, $true, 1, $false, "66-1=", 2) | Out-Null
$d.Content.Find.Execute("3+6=", $true, $false, $false, $false, $false, $true, 1, $false, "96-96=", 2) | Out-Null
$d.Content.Find.Execute("18-7=", $true, $false, $false, $false, $false, $true, 1, $false, "50-30=", 2) | Out-Null
$d.Content.Find.Execute("15+53=", $true, $false, $false, $false, $false, $true, 1, $false, "15+80=", 2) | Out-Null
$d.Content.Find.Execute("23+29=", $true, $false, $false, $false, $false, $true, 1, $false, "33+34=", 2) | Out-Null
$d.Content.Find.Execute("21+7=", $true, $false, $false, $false, $false, $true, 1, $false, "49+32=", 2) | Out-Null
$d.Content.Find.Execute("78-24=", $true, $false, $false, $false, $false, $true, 1, $false, "0+72=", 2) | Out-Null
$d.Content.Find.Execute("45+51=", $true, $false, $false, $false, $false, $true, 1, $false, "22-0=", 2) | Out-Null
